# Auto-generated edit script: updates column F ("想去人数") values
# across all 4 worksheets per the commit diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$updates = @(
    @(2, 1867),
    @(5, 64),
    @(6, 693),
    @(7, 96),
    @(9, 833),
    @(10, 1551),
    @(11, 1242),
    @(12, 1487),
    @(13, 46),
    @(14, 1394),
    @(15, 325),
    @(16, 1646),
    @(18, 796),
    @(19, 1067),
    @(23, 1589),
    @(25, 184),
    @(26, 810),
    @(27, 554),
    @(28, 1158),
    @(29, 0),
    @(30, 1018),
    @(31, 60),
    @(32, 563),
    @(34, 1102),
    @(35, 895),
    @(36, 1102),
    @(37, 55),
    @(38, 221),
    @(39, 61),
    @(40, 857),
    @(41, 1643),
    @(42, 104),
    @(43, 72),
    @(44, 820),
    @(45, 90),
    @(47, 114),
    @(48, 29)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], 6).Value = $u[1]
}

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$updates = @(
    @(5, 180),
    @(6, 1482),
    @(9, 2558),
    @(10, 1204),
    @(12, 721),
    @(13, 238),
    @(18, 449),
    @(21, 309),
    @(22, 0),
    @(29, 195),
    @(31, 45),
    @(34, 176),
    @(41, 134),
    @(42, 60)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], 6).Value = $u[1]
}

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$updates = @(
    @(5, 2859),
    @(6, 4597),
    @(7, 129),
    @(9, 563),
    @(10, 701),
    @(11, 450),
    @(12, 304),
    @(13, 970),
    @(14, 259),
    @(15, 595)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], 6).Value = $u[1]
}

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$updates = @(
    @(2, 1867),
    @(4, 4597),
    @(5, 701),
    @(6, 450),
    @(7, 304),
    @(8, 304),
    @(9, 970),
    @(10, 970),
    @(11, 693),
    @(13, 833),
    @(14, 2558),
    @(15, 1204),
    @(16, 1551),
    @(17, 1242),
    @(18, 1487),
    @(19, 1394),
    @(20, 238),
    @(21, 325),
    @(23, 1646),
    @(24, 796),
    @(25, 1067),
    @(27, 595),
    @(28, 595),
    @(29, 1589),
    @(31, 184),
    @(32, 810),
    @(33, 554),
    @(34, 1158),
    @(35, 309),
    @(36, 1018),
    @(37, 60),
    @(38, 1102),
    @(39, 895),
    @(40, 1102),
    @(42, 221),
    @(43, 857),
    @(45, 1643),
    @(46, 104),
    @(47, 176),
    @(48, 72),
    @(49, 820),
    @(52, 114),
    @(53, 60)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], 6).Value = $u[1]
}
